# Generate Report for Handoff
#
# The "b.md" file has finished a new handoff cycle:
#   - Overview / zh-cn / de-de sheets' status for b.md moves from
#     "Handed back: in sync with en-US" to "Ready for handoff"
#   - new handoff xliff files + timestamps are recorded for zh-cn and de-de
#   - an "Error Detail" note about a stale handback file is attached
#
$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/24c19dd51ec13e0c85b119a16f99ba3f27123333/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3fa21b8ca4d65ce8e846f299f3434ab894edb8ac/e2e/b.md."

# ---------------------------------------------------------------------------
# Overview sheet: b.md is row 3; zh-cn (E) / de-de (F) status + the shared
# "Latest HO Xliff Generate Date" (G) column.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-03 02:41:14"

# ---------------------------------------------------------------------------
# zh-cn sheet: b.md is row 3.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-03 02:41:07"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: b.md is row 3.
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-03 02:41:14"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
